# Auto-generated Excel COM-interop script applying the Shiva_Profits diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the
# affected Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, and
# clears the now-unused HQ-profit cell (column N) on rows that lost their HQ case.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 40000532
$ws.Cells.Item(28, 9).Value = 40000532
$ws.Cells.Item(28, 11).Value = 40000532
$ws.Cells.Item(28, 13).Value = -40000047

# Row 98
$ws.Cells.Item(98, 8).Value = 927.7954999999999
$ws.Cells.Item(98, 9).Value = 1024.64
$ws.Cells.Item(98, 11).Value = 1024.64
$ws.Cells.Item(98, 13).Value = 473.3599999999999

# Row 105
$ws.Cells.Item(105, 8).Value = 88332.664
$ws.Cells.Item(105, 10).Value = 85000
$ws.Cells.Item(105, 12).Value = 85000
$ws.Cells.Item(105, 14).Value = -91988

# Row 113
$ws.Cells.Item(113, 8).Value = 4660.268
$ws.Cells.Item(113, 9).Value = 4112.6665
$ws.Cells.Item(113, 10).Value = 4976.1924
$ws.Cells.Item(113, 11).Value = 4112.6665
$ws.Cells.Item(113, 12).Value = 4976.1924
$ws.Cells.Item(113, 13).Value = -858.6665000000003
$ws.Cells.Item(113, 14).Value = -11484.1924

# Row 122
$ws.Cells.Item(122, 8).Value = 927.7954999999999
$ws.Cells.Item(122, 9).Value = 1024.64
$ws.Cells.Item(122, 11).Value = 3073.92
$ws.Cells.Item(122, 13).Value = -623.9200000000001

# Row 125
$ws.Cells.Item(125, 8).Value = 2865.2856
$ws.Cells.Item(125, 9).Value = 854.6667
$ws.Cells.Item(125, 10).Value = 4373.25
$ws.Cells.Item(125, 11).Value = 7692.0003
$ws.Cells.Item(125, 12).Value = 39359.25
$ws.Cells.Item(125, 13).Value = -5232.0003
$ws.Cells.Item(125, 14).Value = -44279.25

# Row 132
$ws.Cells.Item(132, 8).Value = 6235.229
$ws.Cells.Item(132, 9).Value = 3385.25
$ws.Cells.Item(132, 10).Value = 37585
$ws.Cells.Item(132, 11).Value = 10155.75
$ws.Cells.Item(132, 12).Value = 112755
$ws.Cells.Item(132, 13).Value = -7625.75
$ws.Cells.Item(132, 14).Value = -117815

# Row 138
$ws.Cells.Item(138, 8).Value = 2523.1187
$ws.Cells.Item(138, 9).Value = 1457.2273
$ws.Cells.Item(138, 10).Value = 3156.8918
$ws.Cells.Item(138, 11).Value = 4371.6819
$ws.Cells.Item(138, 12).Value = 9470.6754
$ws.Cells.Item(138, 13).Value = 768.3181000000004
$ws.Cells.Item(138, 14).Value = -19750.6754

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Cells.Item(63, 8).Value = 2912
$ws.Cells.Item(63, 9).Value = 2912
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 2912
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = -2226
$ws.Cells.Item(63, 14).ClearContents()

# Row 66
$ws.Cells.Item(66, 8).Value = 2912
$ws.Cells.Item(66, 9).Value = 2912
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 14560
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = -11128
$ws.Cells.Item(66, 14).Value = -11128

# Row 74
$ws.Cells.Item(74, 8).Value = 2103
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).ClearContents()

# Row 77
$ws.Cells.Item(77, 8).Value = 2103
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).ClearContents()

# Row 122
$ws.Cells.Item(122, 8).Value = 19746.2
$ws.Cells.Item(122, 9).Value = 3059.9167
$ws.Cells.Item(122, 10).Value = 44775.625
$ws.Cells.Item(122, 11).Value = 9179.750100000001
$ws.Cells.Item(122, 12).Value = 134326.875
$ws.Cells.Item(122, 13).Value = -6729.750100000001
$ws.Cells.Item(122, 14).Value = -139226.875

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 2263.6875
$ws.Cells.Item(86, 9).Value = 2570
$ws.Cells.Item(86, 11).Value = 2570
$ws.Cells.Item(86, 13).Value = -1447

# Row 89
$ws.Cells.Item(89, 8).Value = 2263.6875
$ws.Cells.Item(89, 9).Value = 2570
$ws.Cells.Item(89, 11).Value = 12850
$ws.Cells.Item(89, 13).Value = -7234

# Row 105
$ws.Cells.Item(105, 8).Value = 2011.1177
$ws.Cells.Item(105, 9).Value = 2011.1177
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 2011.1177
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -264.1177
$ws.Cells.Item(105, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2912
$ws.Cells.Item(31, 9).Value = 2949.6667
$ws.Cells.Item(31, 11).Value = 2949.6667
$ws.Cells.Item(31, 13).Value = -2654.6667

# Row 34
$ws.Cells.Item(34, 8).Value = 2912
$ws.Cells.Item(34, 9).Value = 2949.6667
$ws.Cells.Item(34, 11).Value = 2949.6667
$ws.Cells.Item(34, 13).Value = -2747.6667

# Row 58
$ws.Cells.Item(58, 8).Value = 1944.091
$ws.Cells.Item(58, 9).Value = 1046.5
$ws.Cells.Item(58, 10).Value = 2457
$ws.Cells.Item(58, 11).Value = 1046.5
$ws.Cells.Item(58, 12).Value = 2457
$ws.Cells.Item(58, 13).Value = -843.5
$ws.Cells.Item(58, 14).Value = -2863

# Row 136
$ws.Cells.Item(136, 8).Value = 1944.091
$ws.Cells.Item(136, 9).Value = 1046.5
$ws.Cells.Item(136, 10).Value = 2457
$ws.Cells.Item(136, 11).Value = 3139.5
$ws.Cells.Item(136, 12).Value = 7371
$ws.Cells.Item(136, 13).Value = -589.5
$ws.Cells.Item(136, 14).Value = -12471

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Cells.Item(9, 8).Value = 59536428
$ws.Cells.Item(9, 10).Value = 125127290
$ws.Cells.Item(9, 12).Value = 375381870
$ws.Cells.Item(9, 14).Value = -375382318

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 6699.6665
$ws.Cells.Item(70, 9).Value = 5824.5
$ws.Cells.Item(70, 10).Value = 8450
$ws.Cells.Item(70, 11).Value = 5824.5
$ws.Cells.Item(70, 12).Value = 8450
$ws.Cells.Item(70, 13).Value = -5554.5
$ws.Cells.Item(70, 14).Value = -8990

# Row 73
$ws.Cells.Item(73, 8).Value = 6699.6665
$ws.Cells.Item(73, 9).Value = 5824.5
$ws.Cells.Item(73, 10).Value = 8450
$ws.Cells.Item(73, 11).Value = 5824.5
$ws.Cells.Item(73, 12).Value = 8450
$ws.Cells.Item(73, 13).Value = -4888.5
$ws.Cells.Item(73, 14).Value = -10322

# Row 102
$ws.Cells.Item(102, 8).Value = 4476.8857
$ws.Cells.Item(102, 9).Value = 4388.5
$ws.Cells.Item(102, 10).Value = 4732.222
$ws.Cells.Item(102, 11).Value = 4388.5
$ws.Cells.Item(102, 12).Value = 4732.222
$ws.Cells.Item(102, 13).Value = -2766.5
$ws.Cells.Item(102, 14).Value = -7976.222

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Cells.Item(55, 8).Value = 433.72223
$ws.Cells.Item(55, 9).Value = 450.5
$ws.Cells.Item(55, 10).Value = 299.5
$ws.Cells.Item(55, 11).Value = 450.5
$ws.Cells.Item(55, 12).Value = 299.5
$ws.Cells.Item(55, 13).Value = -277.5
$ws.Cells.Item(55, 14).Value = -645.5

# Row 82
$ws.Cells.Item(82, 8).Value = 2219.69
$ws.Cells.Item(82, 9).Value = 2252.5474
$ws.Cells.Item(82, 10).Value = 1595.4
$ws.Cells.Item(82, 11).Value = 2252.5474
$ws.Cells.Item(82, 12).Value = 1595.4
$ws.Cells.Item(82, 13).Value = -1891.5474
$ws.Cells.Item(82, 14).Value = -2317.4

# Row 85
$ws.Cells.Item(85, 8).Value = 2219.69
$ws.Cells.Item(85, 9).Value = 2252.5474
$ws.Cells.Item(85, 10).Value = 1595.4
$ws.Cells.Item(85, 11).Value = 2252.5474
$ws.Cells.Item(85, 12).Value = 1595.4
$ws.Cells.Item(85, 13).Value = -1004.5474
$ws.Cells.Item(85, 14).Value = -4091.4

$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Cells.Item(16, 8).Value = 165567.8
$ws.Cells.Item(16, 10).Value = 176959.75
$ws.Cells.Item(16, 12).Value = 176959.75
$ws.Cells.Item(16, 14).Value = -177543.75

# Row 107
$ws.Cells.Item(107, 8).Value = 495.03775
$ws.Cells.Item(107, 9).Value = 422.36365
$ws.Cells.Item(107, 10).Value = 850.3333
$ws.Cells.Item(107, 11).Value = 1267.09095
$ws.Cells.Item(107, 12).Value = 2550.9999
$ws.Cells.Item(107, 13).Value = 652.90905
$ws.Cells.Item(107, 14).Value = -6390.9999
